$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

$rng1 = $ws.Range("B1")
$rng1.Font.Bold = $true
$rng1.Borders.LineStyle = 1
$rng1.HorizontalAlignment = -4108
$rng1.VerticalAlignment = -4160

$rng2 = $ws.Range("A2")
$rng2.Font.Bold = $true
$rng2.Borders.LineStyle = 1
$rng2.HorizontalAlignment = -4108
$rng2.VerticalAlignment = -4160
